$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1888111888111888
$ws.Range("C2").Value = 0.5734265734265734
$ws.Range("J2").Value = 0.006993006993006993
$ws.Range("P2").Value = 0.1223776223776224
$ws.Range("S2").Value = 0.1083916083916084
$ws.Range("B3").Value = 0.005617977528089887
$ws.Range("C3").Value = 0.03370786516853932
$ws.Range("J3").Value = 0.01685393258426966
$ws.Range("P3").Value = 0.7696629213483146
$ws.Range("S3").Value = 0.1741573033707865
$ws.Range("J4").Value = 0.03846153846153846
$ws.Range("P4").Value = 0.6923076923076923
$ws.Range("S4").Value = 0.2692307692307692
$ws.Range("J5").Value = 0.3333333333333333
$ws.Range("P5").Value = 0.6666666666666666
$ws.Range("B6").Value = 0.05116279069767442
$ws.Range("D6").Value = 0.01395348837209302
$ws.Range("F6").Value = 0.06046511627906977
$ws.Range("J6").Value = 0.2883720930232558
$ws.Range("O6").Value = 0.03720930232558139
$ws.Range("Q6").Value = 0.1488372093023256
$ws.Range("R6").Value = 0.08837209302325581
$ws.Range("S6").Value = 0.3116279069767442
$ws.Range("B7").Value = 0.1286549707602339
$ws.Range("D7").Value = 0.01754385964912281
$ws.Range("F7").Value = 0.07602339181286549
$ws.Range("J7").Value = 0.1461988304093567
$ws.Range("O7").Value = 0.02923976608187134
$ws.Range("Q7").Value = 0.1812865497076023
$ws.Range("R7").Value = 0.06432748538011696
$ws.Range("S7").Value = 0.3567251461988304
$ws.Range("B8").Value = 0.08831168831168831
$ws.Range("D8").Value = 0.01558441558441558
$ws.Range("E8").Value = 0.002597402597402597
$ws.Range("F8").Value = 0.04935064935064935
$ws.Range("J8").Value = 0.1064935064935065
$ws.Range("O8").Value = 0.04155844155844156
$ws.Range("Q8").Value = 0.1922077922077922
$ws.Range("R8").Value = 0.08571428571428572
$ws.Range("S8").Value = 0.4181818181818182
$ws.Range("B9").Value = 0.1155555555555556
$ws.Range("D9").Value = 0.02666666666666667
$ws.Range("F9").Value = 0.03555555555555556
$ws.Range("J9").Value = 0.06222222222222222
$ws.Range("O9").Value = 0.03111111111111111
$ws.Range("Q9").Value = 0.1822222222222222
$ws.Range("R9").Value = 0.08444444444444445
$ws.Range("S9").Value = 0.4622222222222222
$ws.Range("B10").Value = 0.1081730769230769
$ws.Range("D10").Value = 0.02724358974358974
$ws.Range("E10").Value = 0.001602564102564103
$ws.Range("F10").Value = 0.07051282051282051
$ws.Range("J10").Value = 0.09695512820512821
$ws.Range("O10").Value = 0.02644230769230769
$ws.Range("Q10").Value = 0.2283653846153846
$ws.Range("R10").Value = 0.07612179487179487
$ws.Range("S10").Value = 0.3645833333333333
$ws.Range("F11").Value = 0.003745318352059925
$ws.Range("G11").Value = 0.1535580524344569
$ws.Range("J11").Value = 0.101123595505618
$ws.Range("K11").Value = 0.2059925093632959
$ws.Range("L11").Value = 0.5280898876404494
$ws.Range("S11").Value = 0.00749063670411985
$ws.Range("G12").Value = 0.7083333333333334
$ws.Range("J12").Value = 0.2083333333333333
$ws.Range("K12").Value = 0.006944444444444444
$ws.Range("L12").Value = 0.03472222222222222
$ws.Range("S12").Value = 0.04166666666666666
$ws.Range("G13").Value = 0.7446808510638298
$ws.Range("J13").Value = 0.2340425531914894
$ws.Range("S13").Value = 0.02127659574468085
$ws.Range("F15").Value = 0.0119047619047619
$ws.Range("H15").Value = 0.126984126984127
$ws.Range("I15").Value = 0.06746031746031746
$ws.Range("J15").Value = 0.3015873015873016
$ws.Range("K15").Value = 0.05952380952380952
$ws.Range("M15").Value = 0.01587301587301587
$ws.Range("N15").Value = 0.003968253968253968
$ws.Range("O15").Value = 0.06349206349206349
$ws.Range("S15").Value = 0.3492063492063492
$ws.Range("F16").Value = 0.01449275362318841
$ws.Range("H16").Value = 0.1449275362318841
$ws.Range("I16").Value = 0.106280193236715
$ws.Range("J16").Value = 0.4106280193236715
$ws.Range("K16").Value = 0.106280193236715
$ws.Range("M16").Value = 0.02415458937198068
$ws.Range("O16").Value = 0.07729468599033816
$ws.Range("S16").Value = 0.1159420289855072
$ws.Range("F17").Value = 0.02380952380952381
$ws.Range("H17").Value = 0.1363636363636364
$ws.Range("I17").Value = 0.1060606060606061
$ws.Range("J17").Value = 0.4069264069264069
$ws.Range("K17").Value = 0.09740259740259741
$ws.Range("M17").Value = 0.01298701298701299
$ws.Range("O17").Value = 0.0367965367965368
$ws.Range("S17").Value = 0.1796536796536796
$ws.Range("F18").Value = 0.01123595505617977
$ws.Range("H18").Value = 0.1741573033707865
$ws.Range("I18").Value = 0.09550561797752809
$ws.Range("J18").Value = 0.4157303370786517
$ws.Range("K18").Value = 0.07865168539325842
$ws.Range("M18").Value = 0.02247191011235955
$ws.Range("O18").Value = 0.07865168539325842
$ws.Range("S18").Value = 0.1235955056179775
$ws.Range("F19").Value = 0.02522935779816514
$ws.Range("H19").Value = 0.1788990825688073
$ws.Range("I19").Value = 0.09403669724770643
$ws.Range("J19").Value = 0.3853211009174312
$ws.Range("K19").Value = 0.08333333333333333
$ws.Range("M19").Value = 0.02599388379204893
$ws.Range("O19").Value = 0.06957186544342507
$ws.Range("S19").Value = 0.1376146788990826
